$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Trends Status" sheet - update Trend Status numbers
# ---------------------------------------------------------------------------
$wsTrends = $wb.Worksheets.Item("Trends Status")

# Rapid Decline
$wsTrends.Range("B2").Value = 0
$wsTrends.Range("C2").Value = 14
$wsTrends.Range("D2").Value = 0
$wsTrends.Range("E2").Value = 19.4

# Decline
$wsTrends.Range("B3").Value = 3
$wsTrends.Range("C3").Value = 18
$wsTrends.Range("D3").Value = 10.7
$wsTrends.Range("E3").Value = 25

# Stable
$wsTrends.Range("B4").Value = 15
$wsTrends.Range("C4").Value = 32
$wsTrends.Range("D4").Value = 53.6
$wsTrends.Range("E4").Value = 44.4

# Increase
$wsTrends.Range("B5").Value = 4
$wsTrends.Range("C5").Value = 2
$wsTrends.Range("D5").Value = 14.3
$wsTrends.Range("E5").Value = 2.8

# Rapid Increase
$wsTrends.Range("B6").Value = 6
$wsTrends.Range("C6").Value = 6
$wsTrends.Range("D6").Value = 21.4
$wsTrends.Range("E6").Value = 8.300000000000001

# Trend Inconclusive
$wsTrends.Range("B7").Value = 51
$wsTrends.Range("C7").Value = 160

# Insufficient Data
$wsTrends.Range("B8").Value = 374
$wsTrends.Range("C8").Value = 221

# ---------------------------------------------------------------------------
# 2. "Priority Status" sheet - update species counts
# ---------------------------------------------------------------------------
$wsPriority = $wb.Worksheets.Item("Priority Status")

$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# 3. "Species qualification" sheet
# ---------------------------------------------------------------------------
$wsQual = $wb.Worksheets.Item("Species qualification")

$wsQual.Range("A2").Value = "SoIB Assessment"
$wsQual.Range("B2").Value = 453

$wsQual.Range("B3").Value = 79
$wsQual.Range("C3").Value = 28

$wsQual.Range("B4").Value = 232
$wsQual.Range("C4").Value = 72

# ---------------------------------------------------------------------------
# 4. Rename "High Priority break-up" -> "Interannual update - High Pri"
#    and update its values
# ---------------------------------------------------------------------------
$wsInter = $wb.Worksheets.Item("High Priority break-up")
$wsInter.Name = "Interannual update - High Pri"

$wsInter.Range("B2").Value = 73
$wsInter.Range("C2").Value = 70.90000000000001
$wsInter.Range("D2").Value = 73
$wsInter.Range("E2").Value = 80.2

$wsInter.Range("B3").Value = 30
$wsInter.Range("C3").Value = 29.1
$wsInter.Range("D3").Value = 18
$wsInter.Range("E3").Value = 19.8

# ---------------------------------------------------------------------------
# 5. Add new sheet "Major update - High Priority " at the end, containing
#    the original (pre-update) "High Priority break-up" data
# ---------------------------------------------------------------------------
$wsMajor = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsMajor.Name = "Major update - High Priority "

$wsMajor.Range("A1").Value = "Break-up"
$wsMajor.Range("B1").Value = "High Species (no.)"
$wsMajor.Range("C1").Value = "High Species (perc.)"
$wsMajor.Range("D1").Value = "New High Species (no.)"
$wsMajor.Range("E1").Value = "New High Species (perc.)"
$wsMajor.Range("A1:E1").Font.Bold = $true
$wsMajor.Range("A1:E1").HorizontalAlignment = -4108

$wsMajor.Range("A2").Value = "Trend New"
$wsMajor.Range("B2").Value = 3
$wsMajor.Range("C2").Value = 15.8
$wsMajor.Range("D2").Value = 3
$wsMajor.Range("E2").Value = 15.8

$wsMajor.Range("A3").Value = "IUCN"
$wsMajor.Range("B3").Value = 16
$wsMajor.Range("C3").Value = 84.2
$wsMajor.Range("D3").Value = 16
$wsMajor.Range("E3").Value = 84.2
